# Inbox test data: Persona/Inbox header row, Customer Success/Sales Orders data row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Persona"
$ws.Range("B1").Value = "Inbox"
$ws.Range("A2").Value = "Customer Success"
$ws.Range("B2").Value = "Sales Orders"

# Widen column A to fit the new "Customer Success" / "Persona" content
# (target best-fit width is 15.7265625 chars; 14.8 is the closest input that
# this engine's column-width rounding resolves to the nearest value)
$ws.Columns("A").ColumnWidth = 14.8

# Move selection to A2 as part of the scrolling/navigation update
$ws.Range("A2").Select()
